$d = $word.ActiveDocument

# Locate "Communications with" (the word to be shortened plus the following
# space+word) so we can work out exact character offsets without hardcoding
# absolute positions.
$findRng = $d.Content
$findRng.Find.Execute("Communications with", $true, $false, $false, $false, $false, `
                       $true, 1, $false, "", 0)

$wordStart = $findRng.Start
$sPos = $wordStart + 13            # offset of the trailing "s" in "Communications"
$splitPos = $wordStart + 14        # offset right after "Communication" (before " with")

# 1) Insert (move) the hidden "_GoBack" bookmark at the point where the edit
#    will happen. Doing this *before* the text edit keeps the paragraph's
#    surrounding runs (e.g. the later ", Education Minor" run) from being
#    merged back together when the text below is changed.
$bmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# 2) Remove the trailing "s" from "Communications" -> "Communication".
$sRange = $d.Range($sPos, $sPos + 1)
$sRange.Text = ""
